$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 89825
$ws.Range("B2").Value = 'Rafaela Câmara'
$ws.Range("C2").Value = 'Engenharia'
$ws.Range("D2").Value = 'Doenca'
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 45099
$ws.Range("G2").Value = 2383.29
$ws.Range("A3").Value = 88008
$ws.Range("B3").Value = 'Lorenzo Freitas'
$ws.Range("C3").Value = 'Vendas'
$ws.Range("D3").Value = 'Viagem de negocios'
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 45085
$ws.Range("G3").Value = 7377.84
$ws.Range("A4").Value = 61751
$ws.Range("B4").Value = 'Vitor Mendes'
$ws.Range("C4").Value = 'P&D'
$ws.Range("D4").Value = 'Viagem de negocios'
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 45105
$ws.Range("G4").Value = 8845.5
$ws.Range("A5").Value = 42322
$ws.Range("B5").Value = 'Rafael Mendonça'
$ws.Range("C5").Value = 'Operacoes'
$ws.Range("D5").Value = 'Outros'
$ws.Range("E5").Value = 6
$ws.Range("F5").Value = 45093
$ws.Range("G5").Value = 8018.74
$ws.Range("A6").Value = 72609
$ws.Range("B6").Value = 'Leandro Peixoto'
$ws.Range("C6").Value = 'Recursos Humanos'
$ws.Range("D6").Value = 'Viagem de negocios'
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 45089
$ws.Range("G6").Value = 4725.29
$ws.Range("A7").Value = 18798
$ws.Range("B7").Value = 'Pedro Miguel Cardoso'
$ws.Range("C7").Value = 'Juridico'
$ws.Range("D7").Value = 'Problemas pessoais'
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 45104
$ws.Range("G7").Value = 5739.41
$ws.Range("A8").Value = 47132
$ws.Range("B8").Value = 'Pedro Miguel Barros'
$ws.Range("C8").Value = 'Operacoes'
$ws.Range("D8").Value = 'Consulta medica'
$ws.Range("E8").Value = 6
$ws.Range("F8").Value = 45079
$ws.Range("G8").Value = 9368.5
$ws.Range("A9").Value = 90454
$ws.Range("B9").Value = 'Yasmin Rocha'
$ws.Range("C9").Value = 'Engenharia'
$ws.Range("D9").Value = 'Doenca'
$ws.Range("E9").Value = 8
$ws.Range("F9").Value = 45106
$ws.Range("G9").Value = 8169.46
$ws.Range("A10").Value = 26240
$ws.Range("B10").Value = 'Sra. Isis Freitas'
$ws.Range("C10").Value = 'P&D'
$ws.Range("D10").Value = 'Outros'
$ws.Range("E10").Value = 4
$ws.Range("F10").Value = 45083
$ws.Range("G10").Value = 2274.16
$ws.Range("A11").Value = 78712
$ws.Range("B11").Value = 'João Guilherme Sampaio'
$ws.Range("C11").Value = 'P&D'
$ws.Range("D11").Value = 'Outros'
$ws.Range("E11").Value = 8
$ws.Range("F11").Value = 45103
$ws.Range("G11").Value = 8422.7
